$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.821.94'
$ws.Range("E2").Value = '  -3.59%  '
$ws.Range("D3").Value = '2.908.41'
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.59'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("E6").Value = '  -5.96%  '
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("D9").Value = '2.907.67'
$ws.Range("E9").Value = '  -4.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  -2.66%  '
$ws.Range("E11").Value = '  -4.71%  '
$ws.Range("E12").Value = '  -3.93%  '
$ws.Range("E13").Value = '  -3.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.42'
$ws.Range("E14").Value = '  -6.51%  '
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").Value = '3.391.71'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("D17").Value = '60.795.00'
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("E18").Value = '  -5.34%  '
$ws.Range("D19").Value = '2.870.35'
$ws.Range("E19").Value = '  -5.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '428.33'
$ws.Range("E20").Value = '  -5.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.55'
$ws.Range("E21").Value = '  -5.05%  '
$ws.Range("E22").Value = '  -2.37%  '
$ws.Range("E23").Value = '  -5.28%  '
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.91'
$ws.Range("E25").Value = '  -3.85%  '
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.91'
$ws.Range("E27").Value = '  -4.14%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.26'
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.47'
$ws.Range("E33").Value = '  -4.00%  '
$ws.Range("D35").Value = '0.0₃0873'
$ws.Range("E35").Value = '  +1.30%  '
$ws.Range("E36").Value = '  -3.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.60'
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.46'
$ws.Range("E39").Value = '  -2.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.126'
$ws.Range("E40").Value = '  -4.44%  '
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.61'
$ws.Range("E42").Value = '  -5.62%  '
$ws.Range("E43").Value = '  -4.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.51'
$ws.Range("E44").Value = '  -5.54%  '
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '377.89'
$ws.Range("E46").Value = '  -3.69%  '
$ws.Range("D47").Value = '2.700.91'
$ws.Range("E47").Value = '  -0.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.30'
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.22'
$ws.Range("E50").Value = '  -4.42%  '
$ws.Range("E51").Value = '  -2.62%  '
